# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices that look numeric (e.g. '68.00', '0.0270'); Excel's Value
# setter auto-converts such strings to numbers and silently drops trailing zeros.
# Prefixing with a literal apostrophe forces text entry (the classic Excel trick),
# then resetting .Style to 'Normal' clears the resulting quote-prefix style so the
# cell format matches the original (unstyled) text cell exactly.

$ws.Range('D2').Value = "'" + '40.260.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.52%  '
$ws.Range('D3').Value = "'" + '2.348.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.14%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'" + '309.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.65%  '
$ws.Range('D6').Value = "'" + '85.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.77%  '
$ws.Range('D7').Value = "'" + '0.531'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.57%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'" + '0.489'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.40%  '
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('D11').Value = "'" + '30.53'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.53%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = "'" + '2.708.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.29%  '
$ws.Range('D14').Value = "'" + '6.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.24%  '
$ws.Range('D15').Value = "'" + '14.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.69%  '
$ws.Range('D16').Value = "'" + '2.347.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.36%  '
$ws.Range('D17').Value = "'" + '0.758'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.49%  '
$ws.Range('D18').Value = "'" + '40.204.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.52%  '
$ws.Range('D19').Value = "'" + '0.0₃0908'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.59%  '
$ws.Range('E20').Value = '  -4.74%  '
$ws.Range('D21').Value = "'" + '68.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.56%  '
$ws.Range('D22').Value = "'" + '10.82'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.68%  '
$ws.Range('D23').Value = "'" + '236.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('E24').Value = '  -6.75%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -6.43%  '
$ws.Range('D27').Value = "'" + '23.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.24%  '
$ws.Range('D28').Value = "'" + '2.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.58%  '
$ws.Range('E29').Value = '  -4.61%  '
$ws.Range('D30').Value = "'" + '35.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.15%  '
$ws.Range('D31').Value = "'" + '151.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.46%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').Value = "'" + '5.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.62%  '
$ws.Range('D34').Value = "'" + '0.0729'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.90%  '
$ws.Range('E35').Value = '  -4.58%  '
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = "'" + '0.101'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.63%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = "'" + '15.96'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.08%  '
$ws.Range('D40').Value = "'" + '1.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.30%  '
$ws.Range('D41').Value = "'" + '3.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.22%  '
$ws.Range('D42').Value = "'" + '2.29'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.29%  '
$ws.Range('D43').Value = "'" + '1.955.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').Value = "'" + '0.0270'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.93%  '
$ws.Range('D45').Value = "'" + '17.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').Value = "'" + '9.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('E47').Value = '  -8.69%  '
$ws.Range('D48').Value = "'" + '2.564.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.35%  '
$ws.Range('D49').Value = "'" + '93.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('D50').Value = "'" + '71.77'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.90%  '
$ws.Range('D51').Value = "'" + '50.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.80%  '
